# Update cryptos list (prices + 1h volume %) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking price strings need a leading quote so Excel stores them
# as text (matching the sheet's existing inlineStr/text cell type) instead
# of silently converting them to the Number type.

$ws.Range('D2').Value = '45.549.10'
$ws.Range('E2').Value = '  +6.69%  '
$ws.Range('D3').Value = '2.387.33'
$ws.Range('E3').Value = '  +4.95%  '
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('D5').Value = "'319.19"
$ws.Range('E5').Value = '  +3.05%  '
$ws.Range('D6').Value = "'111.65"
$ws.Range('E6').Value = '  +8.13%  '
$ws.Range('D7').Value = "'0.638"
$ws.Range('E7').Value = '  +2.74%  '
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('D9').Value = "'0.631"
$ws.Range('E9').Value = '  +5.62%  '
$ws.Range('D10').Value = "'42.10"
$ws.Range('E10').Value = '  +9.04%  '
$ws.Range('E11').Value = '  +3.97%  '
$ws.Range('D12').Value = "'8.68"
$ws.Range('E12').Value = '  +6.09%  '
$ws.Range('E13').Value = '  +5.25%  '
$ws.Range('D14').Value = "'0.109"
$ws.Range('E14').Value = '  +0.87%  '
$ws.Range('E15').Value = '  +5.25%  '
$ws.Range('D16').Value = '2.742.62'
$ws.Range('E16').Value = '  +4.50%  '
$ws.Range('D17').Value = '2.387.12'
$ws.Range('E17').Value = '  +4.93%  '
$ws.Range('D18').Value = '45.542.97'
$ws.Range('E18').Value = '  +7.68%  '
$ws.Range('E19').Value = '  +5.84%  '
$ws.Range('E20').Value = '  +4.16%  '
$ws.Range('D21').Value = "'12.97"
$ws.Range('E21').Value = '  -2.51%  '
$ws.Range('D22').Value = "'75.19"
$ws.Range('E22').Value = '  +3.32%  '
$ws.Range('D23').Value = "'3.56"
$ws.Range('E23').Value = '  +5.00%  '
$ws.Range('D24').Value = "'269.31"
$ws.Range('E24').Value = '  +2.81%  '
$ws.Range('E25').Value = '  +7.51%  '
$ws.Range('E26').Value = '  -0.51%  '
$ws.Range('D27').Value = "'11.35"
$ws.Range('E27').Value = '  +6.90%  '
$ws.Range('D28').Value = "'7.57"
$ws.Range('E28').Value = '  +9.35%  '
$ws.Range('E29').Value = '  +0.37%  '
$ws.Range('E30').Value = '  +3.95%  '
$ws.Range('D31').Value = "'38.68"
$ws.Range('E31').Value = '  +8.69%  '
$ws.Range('D32').Value = "'0.0948"
$ws.Range('E32').Value = '  +11.12%  '
$ws.Range('D33').Value = "'169.92"
$ws.Range('E33').Value = '  +3.47%  '
$ws.Range('D34').Value = "'3.04"
$ws.Range('E34').Value = '  +19.60%  '
$ws.Range('D35').Value = "'0.134"
$ws.Range('E35').Value = '  +3.20%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = "'0.118"
$ws.Range('E36').Value = '  +6.34%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = "'4.88"
$ws.Range('E37').Value = '  +8.77%  '
$ws.Range('D38').Value = "'3.06"
$ws.Range('E38').Value = '  +12.94%  '
$ws.Range('E39').Value = '  +5.59%  '
$ws.Range('D40').Value = "'3.93"
$ws.Range('E40').Value = '  +6.76%  '
$ws.Range('D41').Value = "'1.74"
$ws.Range('E41').Value = '  +12.66%  '
$ws.Range('D42').Value = "'106.51"
$ws.Range('E42').Value = '  +9.01%  '
$ws.Range('D43').Value = "'13.87"
$ws.Range('E43').Value = '  +17.03%  '
$ws.Range('E44').Value = '  +6.89%  '
$ws.Range('D45').Value = "'71.48"
$ws.Range('E45').Value = '  +4.41%  '
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('D47').Value = "'118.31"
$ws.Range('E47').Value = '  +8.10%  '
$ws.Range('D48').Value = "'5.81"
$ws.Range('E48').Value = '  +13.00%  '
$ws.Range('E49').Value = '  +21.15%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').Value = "'9.27"
$ws.Range('E50').Value = '  +7.93%  '
$ws.Range('B51').Value = 'ordi'
$ws.Range('C51').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D51').Value = "'79.13"
$ws.Range('E51').Value = '  +5.79%  '

# Clear the quote-prefix style Excel applied above so the cell keeps the
# workbook's default (unstyled) formatting, same as its neighbours.
foreach ($addr in @('D5','D6','D7','D9','D10','D12','D14','D21','D22','D23','D24','D27','D28','D31','D32','D33','D34','D35','D36','D37','D38','D40','D41','D42','D43','D45','D47','D48','D50','D51')) {
    $ws.Range($addr).Style = 'Normal'
}

